$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: "Status" text used elsewhere via shared string reuse is
#    updated implicitly below (the zh-cn / de-de sheets hold the actual
#    "Status" cells, not the Overview sheet).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (F) / "Latest Handback File"
#    (G) columns now that handback is complete, refresh the handback status
#    text + timestamp, and (re)create the hyperlinks in row order so the
#    relationship ids line up the way a freshly generated report would.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("H2").Value = "2016-03-21 21:02:41"
$wsZh.Range("H3").Value = "2016-03-21 21:02:41"

$wsZh.Range("F2").Value = "67366a8b-6080-41fa-b81b-6cb6d38cd198.md"
$wsZh.Range("G2").Value = "67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.zh-cn.xlf"
$wsZh.Range("F3").Value = "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md"
$wsZh.Range("G3").Value = "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/67366a8b-6080-41fa-b81b-6cb6d38cd198.md", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad52a675e2acda91fb608f587b62c44da0cdbbba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.zh-cn.xlf", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/67366a8b-6080-41fa-b81b-6cb6d38cd198.md", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad52a675e2acda91fb608f587b62c44da0cdbbba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.zh-cn.xlf", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad52a675e2acda91fb608f587b62c44da0cdbbba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.zh-cn.xlf", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad52a675e2acda91fb608f587b62c44da0cdbbba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.zh-cn.xlf", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.zh-cn.xlf")

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, but the handback datetime for this locale
#    is a distinct (later) timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("H2").Value = "2016-03-21 21:02:47"
$wsDe.Range("H3").Value = "2016-03-21 21:02:47"

$wsDe.Range("F2").Value = "67366a8b-6080-41fa-b81b-6cb6d38cd198.md"
$wsDe.Range("G2").Value = "67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.de-de.xlf"
$wsDe.Range("F3").Value = "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md"
$wsDe.Range("G3").Value = "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.de-de.xlf"

$wsDe.Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/67366a8b-6080-41fa-b81b-6cb6d38cd198.md", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/748cc92ec7558fe87a11ccd2d2b73d1c76c071e0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.de-de.xlf", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/67366a8b-6080-41fa-b81b-6cb6d38cd198.md", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/748cc92ec7558fe87a11ccd2d2b73d1c76c071e0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.de-de.xlf", "", "", "67366a8b-6080-41fa-b81b-6cb6d38cd198.4c1798fc505cf4fb24658657e17b4eaffe53105d.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/748cc92ec7558fe87a11ccd2d2b73d1c76c071e0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.de-de.xlf", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/839054eb61e7ed22bc4f925b880f693f177dd9e8/e2e/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/748cc92ec7558fe87a11ccd2d2b73d1c76c071e0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.de-de.xlf", "", "", "c1417f25-d57e-4b55-9a9d-02eb1f70efd5.5f77132c4b14e7c9257815f69be52939f817dab4.de-de.xlf")
